$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed coin data (price/volume rescraped; a couple of rows
# reordered because GateToken jumped in rank). Price (D) and Volume (E)
# columns hold numeric-looking / percent-looking text, so force each
# touched cell to Text format before assigning, keeping them stored as
# text (matching the original inlineStr cells) instead of letting Excel
# auto-convert them to numbers/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '313.75'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2.77%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '35.15'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.121'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.05%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08147'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '3.01%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.121'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.78%'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.155'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.48%'
$ws.Range("B8").Value = 'KuCoinToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '7.958'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.04%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9305'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.74%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1024'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '4.90%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1941'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '5.10%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08994'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '3.79%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03729'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '4.50%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09902'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.36%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001435'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.63%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005830'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.88%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.470'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.26%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.897'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '5.34%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3410'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.83%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1333'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.25%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.106'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.33%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2218'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.23%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04558'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.12%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001250'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.27%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004693'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-3.89%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-3.71%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0004512'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-5.13%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01944'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '4.69%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04872'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '3.00%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007595'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-4.09%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1385'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.81%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.007889'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '1.87%'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-4.17%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006764'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '7.68%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.12%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '196.12'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '291.48%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001705'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-10.42%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.12%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002005'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.12%'
